{"js": "// Replace the 100 arithmetic-problem texts in the single table, in\n// document (row-major) order. Each cell holds one paragraph/run whose\n// text is an expression like \"12+54=\" that gets swapped for a new one.\n// Assigning Table.values with the full grid preserves each cell's\n// existing run formatting (font/size/alignment) while updating the text.\n\nconst newValues = [\n  [\"86-25=\", \"69-31=\", \"2+36=\", \"42-36=\", \"1+88=\"],\n  [\"80-78=\", \"5+33=\", \"68+15=\", \"60-0=\", \"52+37=\"],\n  [\"99-43=\", \"49-13=\", \"0+83=\", \"52+16=\", \"62-14=\"],\n  [\"30+26=\", \"6+23=\", \"43-30=\", \"27-0=\", \"71-40=\"],\n  [\"99-0=\", \"57+21=\", \"16+77=\", \"5+74=\", \"7+83=\"],\n  [\"81-12=\", \"23-22=\", \"84-47=\", \"36-18=\", \"10+36=\"],\n  [\"24+5=\", \"64-44=\", \"59+31=\", \"69-43=\", \"11+59=\"],\n  [\"73-25=\", \"45+19=\", \"47-27=\", \"20-6=\", \"16+4=\"],\n  [\"14+31=\", \"97-77=\", \"55-45=\", \"85-60=\", \"75+24=\"],\n  [\"86-39=\", \"27+47=\", \"48-5=\", \"25-24=\", \"23+48=\"],\n  [\"75-3=\", \"32+4=\", \"94-6=\", \"38+44=\", \"49+13=\"],\n  [\"12+34=\", \"51-22=\", \"11+67=\", \"60-57=\", \"99-58=\"],\n  [\"80+14=\", \"29+35=\", \"10+72=\", \"92-31=\", \"26+28=\"],\n  [\"98-76=\", \"64-0=\", \"61-30=\", \"64-36=\", \"83-82=\"],\n  [\"19+68=\", \"88-6=\", \"80+5=\", \"49+19=\", \"80-67=\"],\n  [\"3+13=\", \"75-67=\", \"33+0=\", \"70-36=\", \"14+6=\"],\n  [\"86-45=\", \"20+6=\", \"38+9=\", \"7+48=\", \"10-0=\"],\n  [\"29+6=\", \"59+7=\", \"36+18=\", \"0+85=\", \"63+13=\"],\n  [\"45-30=\", \"99-7=\", \"73-50=\", \"58+14=\", \"82-37=\"],\n  [\"99-11=\", \"40-31=\", \"16+78=\", \"72+20=\", \"54+36=\"]\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem texts in the single table, in document\n# (row-major) order. Each table cell holds exactly one paragraph/run whose\n# text is an expression like \"12+54=\" that gets swapped for a new one.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$newValues = @(\n    \"86-25=\",\n    \"69-31=\",\n    \"2+36=\",\n    \"42-36=\",\n    \"1+88=\",\n    \"80-78=\",\n    \"5+33=\",\n    \"68+15=\",\n    \"60-0=\",\n    \"52+37=\",\n    \"99-43=\",\n    \"49-13=\",\n    \"0+83=\",\n    \"52+16=\",\n    \"62-14=\",\n    \"30+26=\",\n    \"6+23=\",\n    \"43-30=\",\n    \"27-0=\",\n    \"71-40=\",\n    \"99-0=\",\n    \"57+21=\",\n    \"16+77=\",\n    \"5+74=\",\n    \"7+83=\",\n    \"81-12=\",\n    \"23-22=\",\n    \"84-47=\",\n    \"36-18=\",\n    \"10+36=\",\n    \"24+5=\",\n    \"64-44=\",\n    \"59+31=\",\n    \"69-43=\",\n    \"11+59=\",\n    \"73-25=\",\n    \"45+19=\",\n    \"47-27=\",\n    \"20-6=\",\n    \"16+4=\",\n    \"14+31=\",\n    \"97-77=\",\n    \"55-45=\",\n    \"85-60=\",\n    \"75+24=\",\n    \"86-39=\",\n    \"27+47=\",\n    \"48-5=\",\n    \"25-24=\",\n    \"23+48=\",\n    \"75-3=\",\n    \"32+4=\",\n    \"94-6=\",\n    \"38+44=\",\n    \"49+13=\",\n    \"12+34=\",\n    \"51-22=\",\n    \"11+67=\",\n    \"60-57=\",\n    \"99-58=\",\n    \"80+14=\",\n    \"29+35=\",\n    \"10+72=\",\n    \"92-31=\",\n    \"26+28=\",\n    \"98-76=\",\n    \"64-0=\",\n    \"61-30=\",\n    \"64-36=\",\n    \"83-82=\",\n    \"19+68=\",\n    \"88-6=\",\n    \"80+5=\",\n    \"49+19=\",\n    \"80-67=\",\n    \"3+13=\",\n    \"75-67=\",\n    \"33+0=\",\n    \"70-36=\",\n    \"14+6=\",\n    \"86-45=\",\n    \"20+6=\",\n    \"38+9=\",\n    \"7+48=\",\n    \"10-0=\",\n    \"29+6=\",\n    \"59+7=\",\n    \"36+18=\",\n    \"0+85=\",\n    \"63+13=\",\n    \"45-30=\",\n    \"99-7=\",\n    \"73-50=\",\n    \"58+14=\",\n    \"82-37=\",\n    \"99-11=\",\n    \"40-31=\",\n    \"16+78=\",\n    \"72+20=\",\n    \"54+36=\"\n)\n\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
